$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new value looks numeric, so Excel
# keeps storing them as text (matching the original inline-string cells)
# instead of silently converting them to numbers.
foreach ($addr in @('D4','D5','D6','D7','D9','D10','D11','D13','D15','D17','D19','D21','D22','D23','D25','D27','D29','D30','D31','D32','D33','D34','D36','D37','D38','D39','D44','D45','D46','D47','D48','D49')) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '43.473.74'
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').Value = '2.339.23'
$ws.Range('E3').Value = '  -1.08%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '305.17'
$ws.Range('E5').Value = '  -1.73%  '
$ws.Range('D6').Value = '101.32'
$ws.Range('E6').Value = '  -3.01%  '
$ws.Range('D7').Value = '0.510'
$ws.Range('E7').Value = '  -3.30%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = '0.509'
$ws.Range('E9').Value = '  -2.74%  '
$ws.Range('D10').Value = '35.34'
$ws.Range('E10').Value = '  -2.51%  '
$ws.Range('D11').Value = '0.0799'
$ws.Range('E11').Value = '  -1.87%  '
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('D13').Value = '6.81'
$ws.Range('E13').Value = '  -3.03%  '
$ws.Range('D14').Value = '2.692.31'
$ws.Range('E14').Value = '  -1.40%  '
$ws.Range('D15').Value = '15.72'
$ws.Range('E15').Value = '  +0.01%  '
$ws.Range('D16').Value = '2.339.90'
$ws.Range('E16').Value = '  -0.86%  '
$ws.Range('D17').Value = '0.809'
$ws.Range('E17').Value = '  -0.71%  '
$ws.Range('D18').Value = '43.373.09'
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('D19').Value = '11.90'
$ws.Range('E19').Value = '  -0.89%  '
$ws.Range('D20').Value = '0.0₃0911'
$ws.Range('E20').Value = '  -2.05%  '
$ws.Range('D21').Value = '6.11'
$ws.Range('E21').Value = '  -2.80%  '
$ws.Range('D22').Value = '68.31'
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('D23').Value = '238.03'
$ws.Range('E23').Value = '  -2.32%  '
$ws.Range('E24').Value = '  -3.19%  '
$ws.Range('D25').Value = '2.54'
$ws.Range('E25').Value = '  -3.38%  '
$ws.Range('E26').Value = '  -0.80%  '
$ws.Range('D27').Value = '25.08'
$ws.Range('E27').Value = '  -3.79%  '
$ws.Range('E28').Value = '  -1.76%  '
$ws.Range('D29').Value = '34.71'
$ws.Range('E29').Value = '  -5.24%  '
$ws.Range('D30').Value = '165.85'
$ws.Range('E30').Value = '  +2.09%  '
$ws.Range('D31').Value = '9.25'
$ws.Range('E31').Value = '  -3.91%  '
$ws.Range('D32').Value = '0.998'
$ws.Range('E32').Value = '  -0.12%  '
$ws.Range('D33').Value = '5.08'
$ws.Range('E33').Value = '  -4.28%  '
$ws.Range('D34').Value = '4.57'
$ws.Range('E34').Value = '  -1.43%  '
$ws.Range('E35').Value = '  -4.71%  '
$ws.Range('B36').Value = 'Celestia'
$ws.Range('C36').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D36').Value = '16.92'
$ws.Range('E36').Value = '  -8.05%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = '0.0707'
$ws.Range('E37').Value = '  -4.70%  '
$ws.Range('D38').Value = '2.92'
$ws.Range('E38').Value = '  -6.58%  '
$ws.Range('D39').Value = '1.83'
$ws.Range('E39').Value = '  -6.42%  '
$ws.Range('E40').Value = '  -3.60%  '
$ws.Range('E41').Value = '  -3.19%  '
$ws.Range('E42').Value = '  +0.11%  '
$ws.Range('D43').Value = '1.979.18'
$ws.Range('E43').Value = '  -1.12%  '
$ws.Range('D44').Value = '0.0285'
$ws.Range('E44').Value = '  -2.70%  '
$ws.Range('D45').Value = '18.58'
$ws.Range('E45').Value = '  -6.44%  '
$ws.Range('D46').Value = '9.98'
$ws.Range('E46').Value = '  -4.37%  '
$ws.Range('D47').Value = '2.94'
$ws.Range('E47').Value = '  -5.19%  '
$ws.Range('D48').Value = '56.22'
$ws.Range('E48').Value = '  -3.41%  '
$ws.Range('D49').Value = '4.87'
$ws.Range('E49').Value = '  +3.26%  '
$ws.Range('D50').Value = '2.560.77'
$ws.Range('E50').Value = '  -0.33%  '
$ws.Range('E51').Value = '  -1.80%  '
